# This script inserts one new data row into the "Papa" sheet at row 175,
# pushing the existing rows 175-232 down to 176-233, and populates the
# new row with the appropriate values (mirrors a new weekly observation
# added to the dataset, per the commit message "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 175. This shifts every
# row from 175 downward by one position (175->176, ..., 232->233).
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new record's data.
$ws.Range("A175").Value = 7
$ws.Range("B175").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C175").Value = "Ñuble"
$ws.Range("D175").Value = 44524
$ws.Range("D175").NumberFormat = $ws.Range("D176").NumberFormat
$ws.Range("E175").Value = 16
$ws.Range("F175").Value = 100114001
$ws.Range("G175").Value = "Papa"
$ws.Range("H175").Value = "Patagonia"
$ws.Range("I175").Value = "1a nueva(o)"
$ws.Range("J175").Value = 160
$ws.Range("K175").Value = 9000
$ws.Range("L175").Value = 9500
$ws.Range("M175").Value = 9250
$ws.Range("N175").Value = "$/saco 25 kilos"
$ws.Range("O175").Value = "Región del Maule"
$ws.Range("P175").Value = 370
$ws.Range("Q175").Value = 25
$ws.Range("R175").Value = "Hortaliza"
